$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old sample rows 4-9 (columns A:E), keeping the date-format style on column B
$ws.Range("A4:E9").ClearContents()

# Add a 10th (empty, date-styled) row like the others, by copying B9's format down to B10
$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 1 (header): add new "id" column F
$ws.Range("F1").Value = "id"

# Row 2: value becomes a real number, recipient name updated
$ws.Range("D2").Value = 300.12
$ws.Range("E2").Value = "Claudio Pontes Montarr"

# Row 3: now holds the "Compra de Componentes Eletrônicos" record
$ws.Range("A3").Value = "Compra de Componentes Eletrônicos"
$ws.Range("C3").Value = "Crédito"
$ws.Range("D3").Value = 300.12
$ws.Range("F3").Value = "1593345a-9e7f-449b-b1ab-c9e33a4fea6a"
$ws.Range("E3").Value = "Adele Fonseca"

# Column width adjustments (closest achievable values given engine's width rounding)
$ws.Range("E1").ColumnWidth = 22.666666666666668
$ws.Range("F1").ColumnWidth = 38.666666666666664

# Update the active selection
$ws.Range("E3").Select() | Out-Null
